# NIT-9009945713.xlsx update
# - "Elimina EC anteriores y se agregan nuevos" -> the "Periodo Mora" column
#   (E17:E30) is refreshed so the most recent period (2409) is listed first
#   and the oldest (2308) last, i.e. the previous chronological order is
#   reversed.
# - "se modifica base de datos" -> the "Valor Mora" amount for the first
#   worker row (G16) is corrected from 895000 to 815000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2409","2408","2407","2406","2405","2404","2403","2402","2401","2312","2311","2310","2309","2308")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

$ws.Cells.Item(16, 7).Value = 815000
